# Update "想去人数" (interest count) values in F column across the four
# worksheets of the workbook: 展览 (Exhibitions), 演出 (Performances),
# 本地生活 (Local life), 全部类型 (All types).

$wb = $excel.ActiveWorkbook

# -- Sheet 1: 展览 --
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value2 = 1845
$ws1.Range("F3").Value2 = 408
$ws1.Range("F5").Value2 = 864
$ws1.Range("F6").Value2 = 385
$ws1.Range("F7").Value2 = 748
$ws1.Range("F8").Value2 = 13239
$ws1.Range("F9").Value2 = 13104
$ws1.Range("F10").Value2 = 1003
$ws1.Range("F11").Value2 = 771
$ws1.Range("F15").Value2 = 651
$ws1.Range("F17").Value2 = 59
$ws1.Range("F18").Value2 = 36
$ws1.Range("F19").Value2 = 51
$ws1.Range("F21").Value2 = 214
$ws1.Range("F23").Value2 = 744

# -- Sheet 2: 演出 --
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F9").Value2 = 18

# -- Sheet 3: 本地生活 --
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value2 = 19

# -- Sheet 4: 全部类型 --
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value2 = 1845
$ws4.Range("F4").Value2 = 408
$ws4.Range("F6").Value2 = 864
$ws4.Range("F7").Value2 = 385
$ws4.Range("F9").Value2 = 748
$ws4.Range("F10").Value2 = 13239
$ws4.Range("F11").Value2 = 13104
$ws4.Range("F12").Value2 = 1003
$ws4.Range("F13").Value2 = 771
$ws4.Range("F17").Value2 = 651
$ws4.Range("F21").Value2 = 59
$ws4.Range("F22").Value2 = 36
$ws4.Range("F23").Value2 = 51
$ws4.Range("F27").Value2 = 19
$ws4.Range("F28").Value2 = 214
$ws4.Range("F30").Value2 = 744
$ws4.Range("F33").Value2 = 18
